$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ark1")

$ws.Range("K3").Value = 0
$ws.Range("K4").Value = 0
$ws.Range("K6").Value = 0.5
$ws.Range("K16").Value = 15
$ws.Range("K18").Value = 14.5
$ws.Range("K19").Value = 35
$ws.Range("K21").Value = 9

$co = $ws.ChartObjects(1)
$chart = $co.Chart
$chart.SeriesCollection(2).Values = $ws.Range("D23:K23")
